$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("System State machine")
$ws3.Range("A1:E11").Clear()

$c = $ws3.Range("B1")
$c.Value = "Notes"
$c.Interior.Color = 10086143

$c = $ws3.Range("C1")
$c.Value = "NEW FACE"
$c.Interior.Color = 10086143
$c.Font.Bold = $true

$c = $ws3.Range("D1")
$c.Value = "LIPO INT"
$c.Interior.Color = 10086143
$c.Font.Bold = $true

$c = $ws3.Range("E1")
$c.Value = "TIMER INT"
$c.Interior.Color = 10086143
$c.Font.Bold = $true

$c = $ws3.Range("F1")
$c.Value = "DONE/TimeOut"
$c.Interior.Color = 10086143
$c.Font.Bold = $true

$c = $ws3.Range("A2")
$c.Value = "Initialisation"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("A3")
$c.Value = "STATE_CONFIG"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B3")
$c.Value = "User in config mode, `nAutomaticn face detection disabled"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws3.Range("C3")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("D3")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("E3")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("F3")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A4")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B4")
$c.Value = "Nothing happening until it does"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true

$c = $ws3.Range("C4")
$c.Value = "STATE_CHANGE_TASK"
$c.Interior.Color = 10086143

$c = $ws3.Range("D4")
$c.Value = "STATE_BATTERY_TEST"
$c.Interior.Color = 10086143

$c = $ws3.Range("E4")
$c.Value = "STATE_UPDATE_TASK"
$c.Interior.Color = 10086143

$c = $ws3.Range("F4")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A5")
$c.Value = "STATE_CHANGE_TASK"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B5")
$c.Value = "Set up new Task relating to new face"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true

$c = $ws3.Range("C5")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("D5")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("E5")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("F5")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A6")
$c.Value = "STATE_END_TASK"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B6")
$c.Value = "Dodeca has been placed on STOP Face"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true

$c = $ws3.Range("C6")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("D6")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("E6")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("F6")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A7")
$c.Value = "STATE_UPDATE_TASK"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("C7")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("D7")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("E7")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("F7")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A8")
$c.Value = "STATE_BATTERY_TEST"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B8")
$c.Value = "Test and Display Battery"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true

$c = $ws3.Range("C8")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("D8")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("E8")
$c.Value = "x"
$c.Interior.Color = 10086143

$c = $ws3.Range("F8")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143

$c = $ws3.Range("A9")
$c.Value = "STATE_SLEEP"
$c.Interior.Color = 11389944
$c.Font.Bold = $true

$c = $ws3.Range("B9")
$c.Value = "Dodeca is Sleeping"
$c.Interior.Color = 11389944
$c.NumberFormat = "@"
$c.WrapText = $true

$c = $ws3.Range("C9")
$c.Value = "STATE_CHANGE_TASK"
$c.Interior.Color = 10086143

$c = $ws3.Range("F9")
$c.Value = "STATE_IDLE"
$c.Interior.Color = 10086143
